$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 151
$ws.Range("A151").Value = "Circus approximans"
$ws.Range("B151").Value = "https://www.xeno-canto.org/153004/download"
$ws.Range("C151").Value = "Christchurch, New Zealand"
$ws.Range("D151").Value = "Fernand Deroussen"
$ws.Range("E151").Value = "https://www.xeno-canto.org/contributor/UXGZWVYDFE"
$ws.Range("F151").Value = "assets/misc/cc.png"
$ws.Range("G151").Value = "https://creativecommons.org/licenses/by-nc-nd/3.0/"

# Row 152
$ws.Range("A152").Value = "Circus approximans"
$ws.Range("B152").Value = "https://www.xeno-canto.org/153003/download"
$ws.Range("C152").Value = "Christchurch, New Zealand"
$ws.Range("D152").Value = "Fernand Deroussen"
$ws.Range("E152").Value = "https://www.xeno-canto.org/contributor/UXGZWVYDFE"
$ws.Range("F152").Value = "assets/misc/cc.png"
$ws.Range("G152").Value = "https://creativecommons.org/licenses/by-nc-nd/3.0/"

# Row 153
$ws.Range("A153").Value = "Circus assimilis"
$ws.Range("B153").Value = "https://www.xeno-canto.org/194155/download"
$ws.Range("C153").Value = "Toraut, Sulawesi, Indonesia"
$ws.Range("D153").Value = "Frank Lambert"
$ws.Range("E153").Value = "https://www.xeno-canto.org/contributor/YTUXOCTUEM"
$ws.Range("F153").Value = "assets/misc/cc.png"
$ws.Range("G153").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 154
$ws.Range("A154").Value = "Cisticola exilis"
$ws.Range("B154").Value = "https://www.xeno-canto.org/59528/download"
$ws.Range("C154").Value = "Lake Samsonvale, Queensland, Australia"
$ws.Range("D154").Value = "Peter Woodall"
$ws.Range("E154").Value = "https://www.xeno-canto.org/contributor/SILWLBBIFA"
$ws.Range("F154").Value = "assets/misc/cc.png"
$ws.Range("G154").Value = "https://creativecommons.org/licenses/by-nc-nd/2.5/"

# Row 155
$ws.Range("A155").Value = "Cisticola exilis"
$ws.Range("B155").Value = "https://www.xeno-canto.org/389373/download"
$ws.Range("C155").Value = "Pitt Town Lagoon, New South Wales, Australia"
$ws.Range("D155").Value = "Marc Anderson"
$ws.Range("E155").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F155").Value = "assets/misc/cc.png"
$ws.Range("G155").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 156
$ws.Range("A156").Value = "Climacteris erythrops"
$ws.Range("B156").Value = "https://www.xeno-canto.org/172242/download"
$ws.Range("C156").Value = "Newnes Plateau, New South Wales, Australia"
$ws.Range("D156").Value = "Marc Anderson"
$ws.Range("E156").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F156").Value = "assets/misc/cc.png"
$ws.Range("G156").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 157
$ws.Range("A157").Value = "Climacteris erythrops"
$ws.Range("B157").Value = "https://www.xeno-canto.org/98376/download"
$ws.Range("C157").Value = "Lamington National Park, Queensland, Australia"
$ws.Range("D157").Value = "Patrik Åberg"
$ws.Range("E157").Value = "https://www.xeno-canto.org/contributor/BPSDQEOJWG"
$ws.Range("F157").Value = "assets/misc/cc.png"
$ws.Range("G157").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 158
$ws.Range("A158").Value = "Climacteris picumnus"
$ws.Range("B158").Value = "https://www.xeno-canto.org/343748/download"
$ws.Range("C158").Value = "Cunnamulla, Queensland, Australia"
$ws.Range("D158").Value = "Greg McLachlan"
$ws.Range("E158").Value = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$ws.Range("F158").Value = "assets/misc/cc.png"
$ws.Range("G158").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 159
$ws.Range("A159").Value = "Climacteris picumnus"
$ws.Range("B159").Value = "https://www.xeno-canto.org/334462/download"
$ws.Range("C159").Value = "Gluepot Reserve, South Australia, Australia"
$ws.Range("D159").Value = "Marc Anderson"
$ws.Range("E159").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F159").Value = "assets/misc/cc.png"
$ws.Range("G159").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 160
$ws.Range("A160").Value = "Colluricincla harmonica"
$ws.Range("B160").Value = "https://www.xeno-canto.org/210922/download"
$ws.Range("C160").Value = "Gold Creek Reservoir, Queensland, Australia"
$ws.Range("D160").Value = "Mike Williamson"
$ws.Range("E160").Value = "https://www.xeno-canto.org/contributor/PFQCEGABBH"
$ws.Range("F160").Value = "assets/misc/cc.png"
$ws.Range("G160").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 161
$ws.Range("A161").Value = "Colluricincla harmonica"
$ws.Range("B161").Value = "https://www.xeno-canto.org/434562/download"
$ws.Range("C161").Value = "Mount Moffat, Queensland, Australia"
$ws.Range("D161").Value = "Marc Anderson"
$ws.Range("E161").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F161").Value = "assets/misc/cc.png"
$ws.Range("G161").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 162
$ws.Range("A162").Value = "Colluricincla megarhyncha"
$ws.Range("B162").Value = "https://www.xeno-canto.org/349118/download"
$ws.Range("C162").Value = "Maleny, Queensland, Australia"
$ws.Range("D162").Value = "Marc Anderson"
$ws.Range("E162").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F162").Value = "assets/misc/cc.png"
$ws.Range("G162").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"

# Row 163
$ws.Range("A163").Value = "Columba leucomela"
$ws.Range("B163").Value = "https://www.xeno-canto.org/351946/download"
$ws.Range("C163").Value = "Lamington National Park, Queensland, Australia"
$ws.Range("D163").Value = "Greg McLachlan"
$ws.Range("E163").Value = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$ws.Range("F163").Value = "assets/misc/cc.png"
$ws.Range("G163").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

# Row 164
$ws.Range("A164").Value = "Columba livia"
$ws.Range("B164").Value = "https://www.xeno-canto.org/462661/download"
$ws.Range("C164").Value = "Gmina Hańsk, Lublin Voivodeship, Poland"
$ws.Range("D164").Value = "Stanislas Wroza"
$ws.Range("E164").Value = "https://www.xeno-canto.org/contributor/SDPCHKOHRH"
$ws.Range("F164").Value = "assets/misc/cc.png"
$ws.Range("G164").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"

$ws.Range("B156").Select() | Out-Null
